$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "INCIO SANCHEZ PAOLA KATHERINE",
    "GUEVARA IDROGO DENNIS PERCY",
    "TANTALEAN BUSTAMANTE ESTALIN YOEL",
    "HUAYHUA VALDIVIA LUZ EXMILDA",
    "LOZADA ROJAS LUZ ELENA",
    "DELGADO VASQUEZ FLOR MAGALY",
    "VASQUEZ SILVA ALOIS ADOLF",
    "LINARES PEREZ YANASELY",
    "MEDINA TAPIA ANA YULI",
    "CHAVEZ VILLANUEVA SILVIA JANETH",
    "CAMPOS PEREZ YOVERLY",
    "PEREZ LINARES TATHIANA",
    "SOTO LOZANO LUZDINA",
    "MONDRAGON HERNANDEZ WILMER JUNIOR"
)

$values = @(70, 61, 60, 57, 56, 56, 55, 55, 54, 54, 53, 53, 52, 50)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
